$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5945
$ws.Range("F6").Value = 3000
$ws.Range("F7").Value = 1281
$ws.Range("F8").Value = 443
$ws.Range("F11").Value = 17
$ws.Range("F12").Value = 741
$ws.Range("F13").Value = 292
$ws.Range("F14").Value = 4386
$ws.Range("F15").Value = 4386
$ws.Range("F16").Value = 100
$ws.Range("F18").Value = 123
$ws.Range("F22").Value = 6763
$ws.Range("F24").Value = 107
$ws.Range("F26").Value = 1255
$ws.Range("F27").Value = 6259
$ws.Range("F28").Value = 1641
$ws.Range("F30").Value = 1867
$ws.Range("F31").Value = 6010
$ws.Range("F32").Value = 114
$ws.Range("F34").Value = 101
$ws.Range("F36").Value = 427
$ws.Range("F37").Value = 4524
$ws.Range("F40").Value = 87
$ws.Range("F42").Value = 7
$ws.Range("F43").Value = 2413
$ws.Range("F48").Value = 343
$ws.Range("F49").Value = 2067
$ws.Range("F50").Value = 8
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 195
$ws.Range("F4").Value = 3
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1421
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1421
$ws.Range("F4").Value = 5945
$ws.Range("F5").Value = 3000
$ws.Range("F6").Value = 1281
$ws.Range("F7").Value = 443
$ws.Range("F10").Value = 195
$ws.Range("F12").Value = 292
$ws.Range("F13").Value = 4386
$ws.Range("F14").Value = 4386
$ws.Range("F15").Value = 100
$ws.Range("F17").Value = 123
$ws.Range("F21").Value = 6763
$ws.Range("F23").Value = 107
$ws.Range("F25").Value = 1255
$ws.Range("F27").Value = 6259
$ws.Range("F28").Value = 1641
$ws.Range("F29").Value = 1867
$ws.Range("F31").Value = 6010
$ws.Range("F32").Value = 114
$ws.Range("F35").Value = 101
$ws.Range("F37").Value = 427
$ws.Range("F38").Value = 4525
$ws.Range("F41").Value = 87
$ws.Range("F44").Value = 2413
$ws.Range("F49").Value = 343
